# Auto-generated edit script: applies text replacements via Find/Replace
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-14 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-15 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("581÷6=96, 5", $true, $false, $false, $false, $false, $true, 1, $false, "431÷5=86, 1", 2) | Out-Null
$d.Content.Find.Execute("248÷7=35, 3", $true, $false, $false, $false, $false, $true, 1, $false, "311÷8=38, 7", 2) | Out-Null
$d.Content.Find.Execute("365÷2=182, 1", $true, $false, $false, $false, $false, $true, 1, $false, "239÷3=79, 2", 2) | Out-Null
$d.Content.Find.Execute("749÷9=83, 2", $true, $false, $false, $false, $false, $true, 1, $false, "373÷2=186, 1", 2) | Out-Null
$d.Content.Find.Execute("128÷8=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "699÷3=233, 0", 2) | Out-Null
$d.Content.Find.Execute("680÷6=113, 2", $true, $false, $false, $false, $false, $true, 1, $false, "154÷7=22, 0", 2) | Out-Null
$d.Content.Find.Execute("543÷3=181, 0", $true, $false, $false, $false, $false, $true, 1, $false, "634÷2=317, 0", 2) | Out-Null
$d.Content.Find.Execute("300÷8=37, 4", $true, $false, $false, $false, $false, $true, 1, $false, "680÷7=97, 1", 2) | Out-Null
$d.Content.Find.Execute("467÷6=77, 5", $true, $false, $false, $false, $false, $true, 1, $false, "380÷2=190, 0", 2) | Out-Null
$d.Content.Find.Execute("728÷6=121, 2", $true, $false, $false, $false, $false, $true, 1, $false, "272÷8=34, 0", 2) | Out-Null
$d.Content.Find.Execute("453÷4=113, 1", $true, $false, $false, $false, $false, $true, 1, $false, "470÷4=117, 2", 2) | Out-Null
$d.Content.Find.Execute("234÷7=33, 3", $true, $false, $false, $false, $false, $true, 1, $false, "430÷2=215, 0", 2) | Out-Null
$d.Content.Find.Execute("544÷8=68, 0", $true, $false, $false, $false, $false, $true, 1, $false, "576÷7=82, 2", 2) | Out-Null
$d.Content.Find.Execute("626÷4=156, 2", $true, $false, $false, $false, $false, $true, 1, $false, "949÷9=105, 4", 2) | Out-Null
$d.Content.Find.Execute("345÷5=69, 0", $true, $false, $false, $false, $false, $true, 1, $false, "487÷7=69, 4", 2) | Out-Null
$d.Content.Find.Execute("667÷7=95, 2", $true, $false, $false, $false, $false, $true, 1, $false, "472÷6=78, 4", 2) | Out-Null
$d.Content.Find.Execute("912÷4=228, 0", $true, $false, $false, $false, $false, $true, 1, $false, "746÷2=373, 0", 2) | Out-Null
$d.Content.Find.Execute("919÷9=102, 1", $true, $false, $false, $false, $false, $true, 1, $false, "168÷4=42, 0", 2) | Out-Null
$d.Content.Find.Execute("548÷5=109, 3", $true, $false, $false, $false, $false, $true, 1, $false, "774÷8=96, 6", 2) | Out-Null
$d.Content.Find.Execute("310÷8=38, 6", $true, $false, $false, $false, $false, $true, 1, $false, "572÷6=95, 2", 2) | Out-Null
$d.Content.Find.Execute("920÷4=230, 0", $true, $false, $false, $false, $false, $true, 1, $false, "997÷9=110, 7", 2) | Out-Null
$d.Content.Find.Execute("570÷8=71, 2", $true, $false, $false, $false, $false, $true, 1, $false, "750÷4=187, 2", 2) | Out-Null
$d.Content.Find.Execute("458÷2=229, 0", $true, $false, $false, $false, $false, $true, 1, $false, "503÷4=125, 3", 2) | Out-Null
$d.Content.Find.Execute("195÷4=48, 3", $true, $false, $false, $false, $false, $true, 1, $false, "408÷8=51, 0", 2) | Out-Null
$d.Content.Find.Execute("748÷7=106, 6", $true, $false, $false, $false, $false, $true, 1, $false, "652÷8=81, 4", 2) | Out-Null

Write-Output "done"
